# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# matching the latest export from the logging device (2026-01-28, ~16:30).

$wb = $excel.ActiveWorkbook

function Append-LogRow($ws, $r, $row) {
    # Force text formatting on the whole row first so date-looking /
    # number-looking strings (e.g. "2026-01-28") are stored as literal text,
    # matching the source log export.
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 241-246 (No Motion / Inactive)
# ---------------------------------------------------------------------------
$pirSheet = $wb.Worksheets.Item("PIR")
$pirRows = @(
    ,@("2026-01-28", "16:30:29", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:32", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:34", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:34", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:39", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "16:30:44", "16:00", "Bathroom", "No Motion", "Inactive")
)
$r = 241
foreach ($row in $pirRows) {
    Append-LogRow $pirSheet $r $row
    $r++
}

# ---------------------------------------------------------------------------
# Humidity sheet: rows 235-238 (87.9% / Active)
# ---------------------------------------------------------------------------
$humiditySheet = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    ,@("2026-01-28", "16:30:28", "16:00", "Bathroom", "87.9%", "Active")
    ,@("2026-01-28", "16:30:32", "16:00", "Bathroom", "87.9%", "Active")
    ,@("2026-01-28", "16:30:36", "16:00", "Bathroom", "87.9%", "Active")
    ,@("2026-01-28", "16:30:40", "16:00", "Bathroom", "87.9%", "Active")
)
$r = 235
foreach ($row in $humidityRows) {
    Append-LogRow $humiditySheet $r $row
    $r++
}

# ---------------------------------------------------------------------------
# Temperature sheet: rows 235-239 (22.8C / Active)
# ---------------------------------------------------------------------------
$temperatureSheet = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    ,@("2026-01-28", "16:30:28", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:31", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:33", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:36", "16:00", "Bathroom", "22.8C", "Active")
    ,@("2026-01-28", "16:30:41", "16:00", "Bathroom", "22.8C", "Active")
)
$r = 235
foreach ($row in $temperatureRows) {
    Append-LogRow $temperatureSheet $r $row
    $r++
}
